# Append one new submission row (row 19) to the report sheet, mirroring
# the existing rows' layout: A=notes(blank), B=facilitator, C=quantity,
# D=camp, E=trip type, F=vehicle, G=organization, H=timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(19, 2).Value = "أحمد شريم"

# Quantity "2" must stay text (matches the other rows' numberStoredAsText
# cells) - a leading apostrophe forces text, then copy a plain sibling's
# style so no stray quote-prefix style gets introduced.
$ws.Cells.Item(19, 3).Value = "'2"
$ws.Cells.Item(19, 3).Style = $ws.Cells.Item(18, 3).Style

$ws.Cells.Item(19, 4).Value = "الصمود"
$ws.Cells.Item(19, 5).Value = "الرحلة 2"
$ws.Cells.Item(19, 6).Value = "C2"
$ws.Cells.Item(19, 7).Value = "IDRF"
$ws.Cells.Item(19, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٢٣:٤٥ م"

# Column A is an empty-string cell (not a truly blank/missing cell) in
# every other row, so recreate that: apostrophe forces the cell to exist
# as text, then normalize the style back to the unstyled default.
$ws.Cells.Item(19, 1).Value = "'"
$ws.Cells.Item(19, 1).Style = $ws.Cells.Item(18, 1).Style
